# Add profile update functionality
# - "Accountant" role renamed to "Account man" for Airi Satou, salary bump for
#   Angelica Ramos, "Junior Technical Author" shortened to "Technical Author"
#   (with a salary correction) for Ashton Cox, a profile "updated_at" refresh
#   for every existing employee, a salary correction for Joseph A, and two
#   brand new employees (Sernosh Ulianof, Matt Jinxer) appended to the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated profile fields on existing rows ---
$ws.Range("C1").Value = "Account man"
$ws.Range("G2").Value = 62200
$ws.Range("C3").Value = "Technical Author"
$ws.Range("G3").Value = 3000
$ws.Range("G12").Value = 4000

# --- Refresh the "updated_at" timestamp for every existing employee row ---
$timestamp = "2023-02-19T01:00:19.000000Z"
$ws.Range("F1").Value = $timestamp
$ws.Range("F2").Value = $timestamp
$ws.Range("F3").Value = $timestamp
$ws.Range("F4").Value = $timestamp
$ws.Range("F5").Value = $timestamp
$ws.Range("F6").Value = $timestamp
$ws.Range("F7").Value = $timestamp
$ws.Range("F8").Value = $timestamp
$ws.Range("F9").Value = $timestamp
$ws.Range("F10").Value = $timestamp
$ws.Range("F11").Value = $timestamp
$ws.Range("F12").Value = $timestamp

# --- New employee: Sernosh Ulianof ---
$ws.Range("A13").Value = 13
$ws.Range("B13").Value = "Sernosh Ulianof"
$ws.Range("C13").Value = "Trainer"
$ws.Range("D13").Value = "Chicago"
$ws.Range("E13").Value = 45
$ws.Range("F13").Value = "2023-02-19T14:19:33.000000Z"
$ws.Range("G13").Value = 4000

# --- New employee: Matt Jinxer ---
$ws.Range("A14").Value = 14
$ws.Range("B14").Value = "Matt Jinxer"
$ws.Range("C14").Value = "Integration Specialist"
$ws.Range("D14").Value = "London"
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = "2023-02-19T14:23:18.000000Z"
$ws.Range("G14").Value = 5500
